$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Booking entity block (rows 54-57 before edit) ---
# The sheet was missing an explicit BookingID primary-key row; LocationID and
# EventID were incorrectly marked "Primary Key, Foreign Key, Not Null" instead
# of just "Foreign Key, Not Null". Insert a new row so BookingID gets its own
# row (with trigger/sequence), and LocationID / EventID become plain FK rows.

$ws.Rows("55:55").Insert()

# Row 54: LocationID -> BookingID (Primary Key, Not Null, Unique) + trigger/sequence
$ws.Range("C54").Value2 = "BookingID"
$ws.Range("E54").Value2 = "Primary Key, Not Null, Unique"
$ws.Range("G54").Value2 = "BOOKING_TRIG"
$ws.Range("H54").Value2 = "BOOKING_SEQ"
$ws.Range("K54").Value2 = "Auto Generated Number"

# Row 55 (newly inserted, blank): becomes the LocationID foreign-key row,
# picking up the bordered-cell formatting from I54 so the "Source table"
# column keeps its visual border like its neighbours.
$ws.Range("I54").Copy()
$ws.Range("I55").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("C55").Value2 = "LocationID"
$ws.Range("D55").Value2 = "Number(6)"
$ws.Range("E55").Value2 = "Foreign Key, Not Null"
$ws.Range("I55").Value2 = "None"
$ws.Range("J55").Value2 = "None"
$ws.Range("K55").Value2 = "FK from location database table"

# Row 56 (was row 55, EventID): correct the constraint text only.
$ws.Range("E56").Value2 = "Foreign Key, Not Null"

# Keep the merged label cell (A54:B58) and dimension consistent; Excel grows
# the merge/dimension automatically on row insert, but make sure the
# selection matches the saved workbook state.
$ws.Range("J64").Select()
